$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.509.23'
$ws.Range("E2").Value = '  -1.45%  '

$ws.Range("D3").Value = '1.912.54'
$ws.Range("E3").Value = '  -2.12%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.44'
$ws.Range("E5").Value = '  -1.52%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.0000'
$ws.Range("E6").Value = '  -0.15%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4754'
$ws.Range("E7").Value = '  -2.25%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2843'
$ws.Range("E8").Value = '  -3.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06687'
$ws.Range("E9").Value = '  -4.70%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.79'
$ws.Range("E10").Value = '  -3.63%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '101.20'
$ws.Range("E11").Value = '  -5.74%  '

$ws.Range("D12").Value = '1.920.59'
$ws.Range("E12").Value = '  -1.73%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07678'
$ws.Range("E13").Value = '  -1.01%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.227'
$ws.Range("E14").Value = '  -2.41%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6698'
$ws.Range("E15").Value = '  -4.08%  '

$ws.Range("D16").Value = '30.535.25'
$ws.Range("E16").Value = '  -1.41%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '256.26'
$ws.Range("E17").Value = '  -7.82%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9998'
$ws.Range("E18").Value = '  -0.14%  '

$ws.Range("E19").Value = '  -3.43%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.66'
$ws.Range("E20").Value = '  -4.12%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.400'
$ws.Range("E21").Value = '  -1.38%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.16%  '

$ws.Range("B23").Value = 'Chainlink'
$ws.Range("C23").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.303'
$ws.Range("E23").Value = '  -2.72%  '

$ws.Range("B24").Value = 'Monero'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '168.44'
$ws.Range("E24").Value = '  -0.08%  '

$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.342'
$ws.Range("E25").Value = '  -4.04%  '

$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.99'
$ws.Range("E26").Value = '  -3.35%  '

$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.058'
$ws.Range("E27").Value = '  -4.87%  '

$ws.Range("B28").Value = 'Filecoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.730'
$ws.Range("E28").Value = '  +2.94%  '

$ws.Range("B29").Value = 'Stellar'
$ws.Range("C29").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1007'
$ws.Range("E29").Value = '  -3.59%  '

$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.372'
$ws.Range("E30").Value = '  -2.18%  '

$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.514'
$ws.Range("E31").Value = '  -3.15%  '

$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.256'
$ws.Range("E32").Value = '  -3.05%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04727'
$ws.Range("E33").Value = '  -3.14%  '

$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7279'
$ws.Range("E34").Value = '  -3.11%  '

$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.112'
$ws.Range("E35").Value = '  -4.37%  '

$ws.Range("B36").Value = 'Frax'
$ws.Range("C36").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9990'
$ws.Range("E36").Value = '  -0.12%  '

$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.714'
$ws.Range("E37").Value = '  -0.69%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01914'
$ws.Range("E38").Value = '  -3.94%  '

$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.612'
$ws.Range("E39").Value = '  -2.48%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.224'
$ws.Range("E40").Value = '  -4.60%  '

$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '74.80'
$ws.Range("E41").Value = '  -4.09%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8610'
$ws.Range("E42").Value = '  -3.72%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.950'
$ws.Range("E43").Value = '  -7.34%  '

$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '105.15'
$ws.Range("E44").Value = '  -3.62%  '

$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4241'
$ws.Range("E45").Value = '  -4.26%  '

$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9992'
$ws.Range("E46").Value = '  -0.12%  '

$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '986.31'
$ws.Range("E47").Value = '  -0.45%  '

$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.402'
$ws.Range("E48").Value = '  -5.05%  '

$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1198'
$ws.Range("E49").Value = '  -4.00%  '

$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.81'
$ws.Range("E50").Value = '  -3.04%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.820'
$ws.Range("E51").Value = '  -4.40%  '
